$wb = $excel.ActiveWorkbook

# Sheet references
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The shared string "Ready for handoff" is used by the "Status" columns on all
# three sheets (Overview columns E & F, and column C on the zh-cn / de-de
# sheets). Update every one of those cells to "In Translation" so the shared
# string text is effectively replaced everywhere it appears.
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsDeDe.Range("C2:C3").Value = "In Translation"

# The Status columns got narrower (report regenerated with shorter text), so
# shrink the corresponding columns to match. (12.5 is the character-width
# value that the engine's pixel-snapping rounds to the on-disk width closest
# to the target 13.4101845877511.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
